$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6: week of 9/9 (Tue 9/10/2024), 0.75 hrs, met with Dr. Abuaitah re: logistics/expectations
$ws.Range("A5").Copy()
$ws.Range("A6").PasteSpecial(-4122)
$ws.Range("A6").Value = 45545
$ws.Range("B6").Value = 0.75
$ws.Range("C6").Value = "Met with Dr. Abuaitah to discuss project logistics and expectations"

# Row 7: week of 9/9 (Tue 9/10/2024), 0.5 hrs, met with team to draft contract/roles
$ws.Range("A5").Copy()
$ws.Range("A7").PasteSpecial(-4122)
$ws.Range("A7").Value = 45545
$ws.Range("B7").Value = 0.5
$ws.Range("C7").Value = "Met with team to draft contract and decide roles"

$ws.Range("C10").Select()
